# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.368.54"
$ws.Range("E2").Value = "  -0.61%  "
# Row 3
$ws.Range("D3").Value = "2.641.67"
$ws.Range("E3").Value = "  -1.08%  "
# Row 4
$ws.Range("E4").Value = "  -0.04%  "
# Row 5
$ws.Range("D5").Value = "'587.49"
$ws.Range("E5").Value = "  -2.20%  "
# Row 6
$ws.Range("D6").Value = "'158.02"
$ws.Range("E6").Value = "  +0.69%  "
# Row 7
$ws.Range("D7").Value = "'0.641"
$ws.Range("E7").Value = "  +4.97%  "
# Row 8
$ws.Range("E8").Value = "  -0.01%  "
# Row 9
$ws.Range("D9").Value = "'0.123"
$ws.Range("E9").Value = "  -3.47%  "
# Row 10
$ws.Range("D10").Value = "'5.82"
$ws.Range("E10").Value = "  -0.41%  "
# Row 11
$ws.Range("D11").Value = "'0.393"
$ws.Range("E11").Value = "  -1.26%  "
# Row 12
$ws.Range("E12").Value = "  +0.69%  "
# Row 13
$ws.Range("D13").Value = "'28.82"
$ws.Range("E13").Value = "  -1.30%  "
# Row 14
$ws.Range("D14").Value = "'0.0000188"
$ws.Range("E14").Value = "  -4.10%  "
# Row 15
$ws.Range("D15").Value = "3.115.96"
$ws.Range("E15").Value = "  -1.17%  "
# Row 16
$ws.Range("D16").Value = "65.415.19"
$ws.Range("E16").Value = "  -0.30%  "
# Row 17
$ws.Range("D17").Value = "2.627.93"
$ws.Range("E17").Value = "  -1.05%  "
# Row 18
$ws.Range("D18").Value = "'12.35"
$ws.Range("E18").Value = "  -3.64%  "
# Row 19
$ws.Range("D19").Value = "'4.74"
$ws.Range("E19").Value = "  -1.21%  "
# Row 20
$ws.Range("D20").Value = "'351.95"
$ws.Range("E20").Value = "  +0.09%  "
# Row 21
$ws.Range("D21").Value = "'7.42"
$ws.Range("E21").Value = "  -1.56%  "
# Row 22
$ws.Range("E22").Value = "  -0.02%  "
# Row 23
$ws.Range("D23").Value = "'68.80"
$ws.Range("E23").Value = "  -1.52%  "
# Row 24
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").Value = "'1.75"
$ws.Range("E24").Value = "  +2.90%  "
# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000112"
$ws.Range("E25").Value = "  -0.48%  "
# Row 26
$ws.Range("D26").Value = "'9.47"
$ws.Range("E26").Value = "  -1.48%  "
# Row 27
$ws.Range("E27").Value = "  +1.59%  "
# Row 28
$ws.Range("D28").Value = "'554.53"
$ws.Range("E28").Value = "  +2.87%  "
# Row 29
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.163"
$ws.Range("E29").Value = "  -1.59%  "
# Row 30
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'8.03"
$ws.Range("E30").Value = "  +0.16%  "
# Row 32
$ws.Range("D32").Value = "'2.12"
$ws.Range("E32").Value = "  -0.83%  "
# Row 33
$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  +1.43%  "
# Row 34
$ws.Range("D34").Value = "'6.54"
$ws.Range("E34").Value = "  +0.91%  "
# Row 35
$ws.Range("D35").Value = "'5.43"
$ws.Range("E35").Value = "  -1.20%  "
# Row 36
$ws.Range("D36").Value = "'0.417"
$ws.Range("E36").Value = "  -1.19%  "
# Row 37
$ws.Range("D37").Value = "'20.32"
$ws.Range("E37").Value = "  -0.99%  "
# Row 38
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.03%  "
# Row 39
$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +0.99%  "
# Row 40
$ws.Range("D40").Value = "'152.20"
$ws.Range("E40").Value = "  -3.55%  "
# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'160.48"
$ws.Range("E41").Value = "  -2.17%  "
# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.44"
$ws.Range("E42").Value = "  +5.28%  "
# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'4.06"
$ws.Range("E43").Value = "  -1.11%  "
# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0609"
$ws.Range("E44").Value = "  +0.12%  "
# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'22.99"
$ws.Range("E45").Value = "  +1.31%  "
# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.640"
$ws.Range("E46").Value = "  +0.03%  "
# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0255"
$ws.Range("E47").Value = "  -0.94%  "
# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.102"
$ws.Range("E48").Value = "  +3.03%  "
# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'19.59"
$ws.Range("E49").Value = "  -2.30%  "
# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0241"
$ws.Range("E50").Value = "  -6.71%  "
# Row 51
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'0.795"
$ws.Range("E51").Value = "  -2.51%  "
